$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Set the whole B2:D9 block to 0 first
$ws.Range("B2:D9").Value = 0

# Then apply the specific non-zero overrides from the diff
$ws.Range("C3").Value = 0.634706448907856
$ws.Range("C7").Value = -0.6304474247132246
